$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Figures")

$data = @()
$data += ,@("FIG-01", "Developed markets prioritize 'Value over Volume' in pet ownership trends.", "Chart", "S089", "Figure 1", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-02", "European demographics show a structural dominance of the feline segment.", "Chart", "S089", "Figure 2", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-03", "Feline segment growth outpaces canine due to urbanization constraints.", "Chart", "S089", "Figure 3", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-04", "APAC emerging as the primary volume engine for the next decade.", "Chart", "S089", "Figure 4", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-05", "Probiotics volume share reflects the dominance of poultry and swine sectors.", "Chart", "S089", "Figure 5", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-06", "Disease outbreaks act as catalysts for immune-modulating additive demand.", "Chart", "S089", "Figure 6", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-07", "Regulatory pressures drive structural contraction in Western swine herds.", "Chart", "S089", "Figure 7", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-08", "Western de-ruminization shifts value from herd size to efficiency-per-head.", "Chart", "S089", "Figure 8", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-09", "Global protein production shifts favor poultry and aquaculture over ruminants.", "Chart", "S089", "Figure 9", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-11", "The 'Blue Transformation' drives industrialization and functional additive needs in aqua.", "Chart", "S089", "Figure 11", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-11", "Palatability and format dictate compliance, which is a key efficacy driver.", "Chart", "S089", "Figure 11", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-12", "Nutraceuticals have captured a dominant share of the preventive care wallet.", "Chart", "S089", "Figure 12", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-13", "High-spending households drive the majority of revenue in the pet wellness category.", "Chart", "S089", "Figure 13", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-14", "Purchasing psychology is driven more by the fear of loss than aspirational health.", "Chart", "S089", "Figure 14", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-15", "Category premiumization evolves from generic ingredients to IP-backed solutions.", "Chart", "S089", "Figure 15", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-16", "Determining validity of 'Pre-Senior' segment expands Customer Lifetime Value.", "Chart", "S089", "Figure 16", "Part II", "Final", "2026-02-08", "")
$data += ,@("FIG-18", "Channel economics favor omnichannel dominance but penalize pure-play DTC.", "Chart", "S089", "Figure 18", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-19", "Pricing power erodes as products move from IP owners to generic retailers.", "Chart", "S089", "Figure 19", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-20", "Strategic categorization defines the risk-reward profile for market entrants.", "Chart", "S089", "Figure 20", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-21", "Incumbents utilize nutraceuticals as a low-cost acquisition funnel for future banking.", "Chart", "S089", "Figure 21", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-33", "Figure 33: Margin capture shifts upstream to IP holders and downstream to specialized CDMOs.", "Chart", "S089", "Figure 33", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_0_1_Innovation_Matrix", "R&D intensity correlates directly with EBITDA margin expansion/premium.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_0_2_Market_Bifurcation", "Structural bifurcation splits the market into Emotional (Pet) and ROI (Livestock) economies.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_10_Matrix", "Advanced delivery formats ensure bioavailability and maximize compliance.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_11_Matrix", "Sustainability metrics like methane reduction are becoming non-negotiable procurement specs.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_1_Matrix", "Efficacy levels in Mobility define market positioning and pricing power.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_2_Matrix", "Gut Health strategy shifts from generic digestion to precision microbiome modulation.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_3_Matrix", "Immunity solutions build biological resilience for the post-antibiotic era.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_4_Matrix", "Cognitive support monetizes the Silver Economy via neuro-preservation.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_5_Matrix", "Non-sedative anxiolysis replaces pharmacological interventions in behavior management.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_6_Matrix", "Enzymes and yeast cultures drive Feed Conversion Ratios (FCR) in livestock.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_II_9_Matrix", "Gene-expression data constructs a defensible moat of validation around ingredients.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_I_3_Regulatory_Matrix", "Regulatory pathways dictate unit economics and allowable claims.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_MA_Matrix", "M&A valuation matrix favors de-risked assets over internal R&D development.", "Chart", "S089", "", "Part III", "Final", "2026-02-08", "")
$data += ,@("FIG-Figure_TAM_Reconciliation", "Market reconciliation excludes commodities to define the investable high-value universe.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Regulatory_Timeline_RightArrow", "Regulatory Timeline: The Push for Alternatives", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")
$data += ,@("FIG-Table_US_vs_EU", "Regulatory Divergence creates structural barriers to entry for non-compliant actors.", "Chart", "S089", "", "Part I", "Final", "2026-02-08", "")

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $colIdx = $c + 1
        $val = $rowVals[$c]
        $cell = $ws.Cells.Item($r, $colIdx)
        if ($val -ne "") {
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        }
    }
}
